$d = $word.ActiveDocument

# 1. Title paragraph text change
$d.Paragraphs.Item(1).Range.Text = "ΕΚΘΕΣΗ ΕΝΟΡΚΗΣ ΕΞΕΤΑΣΗΣ"

# 2. Intro paragraph (paragraph 2) full text replace
$d.Paragraphs.Item(2).Range.Text = "    Στην {{place}} σήμερα την {{date_num}} του μήνα {{month}} του έτους {{year}} ημέρα εβδομάδας {{day}} και ώρα {{hour}} ενώπιον εμού του {{first_officer}} του {{policeStation}} Θεσσαλονίκης, παρισταμένου και της {{sec_officer}} της ιδίας υπηρεσίας, που προσλήφθηκε ως Β' Ανακριτικός Υπάλληλος, εμφανίστηκε ο κατωτέρω μάρτυρας, ο οποίος αφού ρωτήθηκε για την ταυτότητα του κ.λ.π. απάντησε ότι ονομάζεται: {{surname}} {{name}} του {{fathername}} και της {{mothername}} γεν. {{dateOfBirth}} στη {{placeOfBirth}} κατ. {{address}}, αριθμός τηλεφώνου {{tel}}, ηλεκτρονικό ταχυδρομείο {{email}}, κάτοχος του υπ αριθμόν {{DAT}} που εκδόθηκε την {{issued}} από {{place_issued}} Α.Φ.Μ : {{afm}}, Δ.Ο.Υ : {{doy}}"

# 3. Closing paragraph (paragraph 3) replaced with oath text
$d.Paragraphs.Item(3).Range.Text = "   Έπειτα ο εξεταζόμενος έδωσε τον προβλεπόμενο από τα άρθρα 219 και 220 παρ. 1  του Κώδικα Ποινικής Δικονομίας όρκο, ως ακολούθως: «Δηλώνω, επικαλούμενος την τιμή και την συνείδηση μου, ότι θα πω όλη την αλήθεια και μόνο την αλήθεια, χωρίς να προσθέσω ούτε να αποκρύψω τίποτα», και στην συνέχεια εξετάσθηκε ως εξής:."

# 4. Insert three new paragraphs after paragraph 3 (oath), before the signature paragraph
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()
$d.Paragraphs.Item(4).Range.Text = "ΕΡΩΤΗΣΗ: Ρωτήθηκε σχετικά:"

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$d.Paragraphs.Item(5).Range.Text = "ΑΠΟΚΡΙΣΗ: Στις {{dateOfCrime}} και περί ώρα {{hourOfCrime}} στη {{placeOfCrime}}, {{stateOfVictim}}, o {{surnamePerperator}} {{namePerperator}} του {{fathernamePerperator}} και της {{mothernamePerperator}}γεν. {{dateOfBirthPerperator}} στη {{placeOfBirthPerperator}} κατ. {{addressPerperator}},αριθμός τηλεφώνου {{telPreperator}}, ηλεκτρονικό  ταχυδρομείου{{emailPreperator}},κάτοχος του υπ αριθμόν {{DATperperator}} που εκδόθηκε την {{issuedPerperator}} από {{place_issuedPerperator}}  Α.Φ.Μ : {{afmPreperator}}, Δ.Ο.Υ : {{doyPrep}}  {{whatHappened}},{{howHappened}}.{{whyHappened}}.{{add_something}}. {{forensicExam}}.{{prosecution}}.Τίποτε άλλο δεν έχω να προσθέσω και υπογράφω,"

$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$d.Paragraphs.Item(6).Range.Text = "Η παρούσα έκθεση άρχισε να συντάσσεται την {{hour}} ώρα και περαιώθηκε την {{hourOfReportFinished}} ώρα. Για πίστωση συντάχθηκε η παρούσα έκθεση η οποία αφούαναγνώσθηκε και βεβαιώθηκε, υπογράφεται ως ακολούθως:"

Write-Output "Paragraph count:"
Write-Output $d.Paragraphs.Count
